# #5: cash & deposit done
# Restructure the "存款" (deposit) worksheet (sheet index 3) from the old
# ad-hoc 7-column layout (bank, deposit_type, currency, owner, <blank qty>, total)
# into the same full metadata layout used by the other asset sheets:
# bank, deposit_type, currency, owner, total, property_category, category,
# date, legislator_name, legislator_id, source_file, index.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---- Header row (row 1) ----
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Re-apply the bold/bordered/centered header formatting across the whole
# header row so the newly-used F1:M1 cells converge back onto the same
# (already-present) header style as B1:E1 instead of staying unformatted.
$headerRow = $ws.Range("B1:M1")
$headerRow.Font.Bold = $true
$headerRow.HorizontalAlignment = -4108
$headerRow.VerticalAlignment = -4160
$headerRow.Borders.LineStyle = 1

# "date" column (I2:I8) holds the literal text "2012-05-01" in every data
# row, not a real Excel date. Force text storage (else Excel auto-converts
# the ISO-looking string into a date serial) and then drop back to the
# default "Normal" style so the cell doesn't keep a stray Text number
# format / style index that the original file never had.
$dateCol = $ws.Range("I2:I8")
$dateCol.NumberFormat = "@"

# ---- Data rows (rows 2-8) ----
# row 2 (index 83): 安泰商業銀行民權分行 / 活期儲蓄存款 / 新臺幣 / 11613676
$ws.Range("B2").Value = "安泰商業銀行民權分行"
$ws.Range("C2").Value = "活期儲蓄存款"
$ws.Range("D2").Value = "新臺幣"
$ws.Range("E2").Value = "羅淑蕾"
$ws.Range("F2").Value = 11613676
$ws.Range("G2").Value = "deposit"
$ws.Range("H2").Value = "normal"
$ws.Range("I2").Value = "2012-05-01"
$ws.Range("J2").Value = "羅淑蕾"
$ws.Range("K2").Value = 1638
$ws.Range("L2").Value = "tmpe6421"
$ws.Range("M2").Value = 83

# row 3 (index 84): 安泰商業銀行 / 支票存款 / 新臺幣 / 161
$ws.Range("B3").Value = "安泰商業銀行"
$ws.Range("C3").Value = "支票存款"
$ws.Range("D3").Value = "新臺幣"
$ws.Range("E3").Value = "羅淑蕾"
$ws.Range("F3").Value = 161
$ws.Range("G3").Value = "deposit"
$ws.Range("H3").Value = "normal"
$ws.Range("I3").Value = "2012-05-01"
$ws.Range("J3").Value = "羅淑蕾"
$ws.Range("K3").Value = 1638
$ws.Range("L3").Value = "tmpe6421"
$ws.Range("M3").Value = 84

# row 4 (index 85): 華泰商業銀行 / 活期儲蓄存款 / 新臺幣 / 25246
$ws.Range("B4").Value = "華泰商業銀行"
$ws.Range("C4").Value = "活期儲蓄存款"
$ws.Range("D4").Value = "新臺幣"
$ws.Range("E4").Value = "羅淑蕾"
$ws.Range("F4").Value = 25246
$ws.Range("G4").Value = "deposit"
$ws.Range("H4").Value = "normal"
$ws.Range("I4").Value = "2012-05-01"
$ws.Range("J4").Value = "羅淑蕾"
$ws.Range("K4").Value = 1638
$ws.Range("L4").Value = "tmpe6421"
$ws.Range("M4").Value = 85

# row 5 (index 86): 華泰商業銀行 / 支票存款 / 新臺幣 / 59831
$ws.Range("B5").Value = "華泰商業銀行"
$ws.Range("C5").Value = "支票存款"
$ws.Range("D5").Value = "新臺幣"
$ws.Range("E5").Value = "羅淑蕾"
$ws.Range("F5").Value = 59831
$ws.Range("G5").Value = "deposit"
$ws.Range("H5").Value = "normal"
$ws.Range("I5").Value = "2012-05-01"
$ws.Range("J5").Value = "羅淑蕾"
$ws.Range("K5").Value = 1638
$ws.Range("L5").Value = "tmpe6421"
$ws.Range("M5").Value = 86

# row 6 (index 87): 中國銀行 / 綜合存款 / 人民幣 / 7879500
$ws.Range("B6").Value = "中國銀行"
$ws.Range("C6").Value = "綜合存款"
$ws.Range("D6").Value = "人民幣"
$ws.Range("E6").Value = "羅淑蕾"
$ws.Range("F6").Value = 7879500
$ws.Range("G6").Value = "deposit"
$ws.Range("H6").Value = "normal"
$ws.Range("I6").Value = "2012-05-01"
$ws.Range("J6").Value = "羅淑蕾"
$ws.Range("K6").Value = 1638
$ws.Range("L6").Value = "tmpe6421"
$ws.Range("M6").Value = 87

# row 7 (index 88): 美商美國銀行 / 綜合存款 / 美金 / 30533705.8
$ws.Range("B7").Value = "美商美國銀行"
$ws.Range("C7").Value = "綜合存款"
$ws.Range("D7").Value = "美金"
$ws.Range("E7").Value = "羅淑蕾"
$ws.Range("F7").Value = 30533705.8
$ws.Range("G7").Value = "deposit"
$ws.Range("H7").Value = "normal"
$ws.Range("I7").Value = "2012-05-01"
$ws.Range("J7").Value = "羅淑蕾"
$ws.Range("K7").Value = 1638
$ws.Range("L7").Value = "tmpe6421"
$ws.Range("M7").Value = 88

# row 8 (index 89): 永豐商業銀行敦南分行 / 活期儲蓄存款 / 新臺幣 / 1593775
$ws.Range("B8").Value = "永豐商業銀行敦南分行"
$ws.Range("C8").Value = "活期儲蓄存款"
$ws.Range("D8").Value = "新臺幣"
$ws.Range("E8").Value = "羅淑蕾"
$ws.Range("F8").Value = 1593775
$ws.Range("G8").Value = "deposit"
$ws.Range("H8").Value = "normal"
$ws.Range("I8").Value = "2012-05-01"
$ws.Range("J8").Value = "羅淑蕾"
$ws.Range("K8").Value = 1638
$ws.Range("L8").Value = "tmpe6421"
$ws.Range("M8").Value = 89

# Strip the temporary Text number format back off I2:I8 so these cells end
# up with the same (absent / default) style as every other newly-added
# data cell in columns G-M, matching the un-styled "total"-column look of
# the rest of the sheet.
$dateCol.Style = "Normal"
